$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.028.25'
$ws.Range('E2').Value = '  -2.00%  '
$ws.Range('D3').Value = '3.581.91'
$ws.Range('E3').Value = '  -3.14%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.52'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -5.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '191.96'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.82%  '
$ws.Range('D7').Value = '3.579.01'
$ws.Range('E7').Value = '  -3.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.615'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.61%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.679'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -5.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '55.69'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -6.86%  '
$ws.Range('E12').Value = '  -6.19%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000270'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -5.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.86'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -4.94%  '
$ws.Range('D15').Value = '4.154.81'
$ws.Range('E15').Value = '  -2.83%  '
$ws.Range('D16').Value = '3.582.95'
$ws.Range('E16').Value = '  -2.80%  '
$ws.Range('E17').Value = '  -1.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.36'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -5.04%  '
$ws.Range('D19').Value = '67.016.22'
$ws.Range('E19').Value = '  -1.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.19'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -4.49%  '
$ws.Range('E21').Value = '  -6.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '400.12'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.65%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.18'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -8.59%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '86.00'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.40'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.95%  '
$ws.Range('E26').Value = '  -3.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.45'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.82%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.09'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.62'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.48%  '
$ws.Range('E30').Value = '  -6.74%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.67'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('E32').Value = '  -3.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '636.59'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.14%  '
$ws.Range('E34').Value = '  -4.32%  '
$ws.Range('E35').Value = '  -5.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '63.95'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -5.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '42.16'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -11.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.401'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.21%  '
$ws.Range('E39').Value = '  +0.31%  '
$ws.Range('D40').Value = '0.0₃0768'
$ws.Range('E40').Value = '  -6.83%  '
$ws.Range('D41').Value = '3.199.61'
$ws.Range('E41').Value = '  +11.04%  '
$ws.Range('E42').Value = '  -3.27%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.70'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.64%  '
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.97'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0415'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -6.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.11'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.95%  '
$ws.Range('E48').Value = '  -6.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '141.90'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.59'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.72%  '
$ws.Range('E51').Value = '  -6.62%  '
